# Applies "updating with Paula's interpretation" edits to the Black Jack
# Simulation write-up: rewording of a couple of bullets, and a large
# expansion of the final bullet into a walk-through of the simulation
# sequence with three worked examples.
#
# NOTE: this interpreter mis-parses calls of the form
#   Some-Function $comObj ("literal " + "text")
# (a parenthesised expression as a later argument), so every text value
# below is assigned to a plain variable first and that variable is
# passed to the helper functions - never an inline parenthesised
# expression.

$d = $word.ActiveDocument
$RIGHT_QUOTE = [char]0x2019

# Replace the full text of a paragraph. Assigning straight to
# Range.Text only ever touches the paragraph's first run, so instead we
# delete the paragraph's content (but keep its trailing mark) and then
# insert the replacement text into the now-empty range.
function Set-ParaText($para, [string]$text) {
    $s = $para.Range.Start
    $e = $para.Range.End - 1
    if ($e -gt $s) {
        $d.Range($s, $e).Delete() | Out-Null
    }
    $d.Range($s, $s).InsertAfter($text)
}

# Insert a blank ListParagraph spacer followed by a new numbered bullet
# (list id 2) containing $text, both right after $afterPara. Both new
# paragraphs are created (inheriting the numbered pPr from $afterPara)
# before the spacer's numbering is stripped, so that stripping the
# spacer's numbering doesn't also strip the bullet's (InsertParagraphAfter
# always clones pPr off of whatever paragraph currently precedes the
# insertion point). Returns the new bullet paragraph.
function Add-SpacerAndBullet($afterPara, [string]$text) {
    $afterPara.Range.InsertParagraphAfter() | Out-Null
    $spacer = $afterPara.Next()
    $spacer.Range.InsertParagraphAfter() | Out-Null
    $bullet = $spacer.Next()
    $spacer.Range.ListFormat.RemoveNumbers() | Out-Null
    Set-ParaText $bullet $text
    return $bullet
}

# ---------------------------------------------------------------------
# 1) "The simulation tests the idea of the player to stay when he has
#    reached 12 points." -> "...staying when he reaches 12 points."
# ---------------------------------------------------------------------
$rng1 = $d.Content
$old1 = "to stay when he has reached 12 points"
$new1 = "staying when he reaches 12 points"
$rng1.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Rewrite the "So if a player hits..." bullet with Paula's wording,
#    and fold the blank separator paragraph that used to follow it back
#    into the list (the new wording doesn't need a spacer there).
# ---------------------------------------------------------------------
$rng2 = $d.Content
$needle2 = "So if a player hits after he reaches 12 there is a chance for the player to receive a 10 point card and bust (lose)"
$rng2.Find.Execute($needle2) | Out-Null
$bullet2 = $rng2.Paragraphs(1)
$text2 = "When a player" + $RIGHT_QUOTE + "s hand is 12, if he hits there is a chance the play will receive a 10  point card and bust.  The purpose of this simulation is to lower the chance of busting."
Set-ParaText $bullet2 $text2

$spacer2 = $bullet2.Next()
$spacer2text = $spacer2.Range.Text.Trim()
if ($spacer2text.Length -eq 0) {
    $spacer2.Range.Delete() | Out-Null
}

# ---------------------------------------------------------------------
# 3) Rewrite the "To increase the chances..." bullet.
# ---------------------------------------------------------------------
$rng3 = $d.Content
$needle3 = "To increase the chances"
$rng3.Find.Execute($needle3) | Out-Null
$bullet3 = $rng3.Paragraphs(1)
$text3 = "If a player" + $RIGHT_QUOTE + "s hand is 11, he can hit without busting because the highest non-interchangeable card is valued at 10.  This hand also has a higher possibility of getting blackjack versus busting."
Set-ParaText $bullet3 $text3

# ---------------------------------------------------------------------
# 4) Expand the final "Through the simulation, " bullet with the rest of
#    that sentence, then append a long walk-through (sequence + three
#    worked examples + trailing bullet), each separated by a blank
#    ListParagraph spacer, matching the surrounding bullets' style.
# ---------------------------------------------------------------------
$rng4 = $d.Content
$needle4 = "Through the simulation"
$rng4.Find.Execute($needle4) | Out-Null
$bullet4 = $rng4.Paragraphs(1)
$tail4 = "there are four decks of cards (52 * 4 = 208) and set the value of Jack, Queen, King at 10 points and Ace is valued at 11 or 1 and it is shuffled randomly."
$bullet4.Range.InsertAfter($tail4)

$seqText = "The sequence starts with the player and dealer receiving two cards that will be less than or equal to 11.  At this point the player can ask for another card. "
$ex1Text = "Example1: Player gets 3 and 4 and ask for another card and gets a jack (10 points) then the total is 17 and stays to give a chance for the dealer to bust."
$ex2Text = "Example2: player gets 3 and 4 and receives another 4 equals to 11 at this point the player can still ask for another card and receives a queen at this point the player hits blackjack."
$ex3Text = "Example3: player gets 5 and 7 it equals to 12 at this point if he asks for another and it" + $RIGHT_QUOTE + "s a 10 point card he can bust (lose). So the plan is to stay at 12 and chance that the dealer will bust."
$lastText = "  "

$cur = $bullet4
$cur = Add-SpacerAndBullet $cur $seqText
$cur = Add-SpacerAndBullet $cur $ex1Text
$cur = Add-SpacerAndBullet $cur $ex2Text
$cur = Add-SpacerAndBullet $cur $ex3Text
$cur = Add-SpacerAndBullet $cur $lastText

Write-Host "Final paragraph count: $($d.Paragraphs.Count)"
